$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in B2/B3
$ws.Range("B2").Value = 0.2
$ws.Range("B3").Value = 0.6926431322749765

# Rows 10-14: labels shift up by one (A14 moves to A10), values updated accordingly
$ws.Range("A10").Value = "Overall 70% Redundancy"
$ws.Range("B10").Value = 0.9111111111111111

$ws.Range("A11").Value = "Disparity Generators"
$ws.Range("B11").Value = 0.6778822900377038

$ws.Range("A12").Value = "Disparity Load"
$ws.Range("B12").Value = 0.5288308596787999

$ws.Range("A13").Value = "Disparity Trafo"
$ws.Range("B13").Value = 0

$ws.Range("A14").Value = "Disparity Lines"
$ws.Range("B14").Value = 0.8101873098292478
